$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab/workbook.xml <sheet name="...">)
$ws.Name = "GossF"

# 2. Tiny value correction on G13 (last-bit precision tweak)
$ws.Range("G13").Value = 0.9951615938329431

# 3. Append new data row 16, mirroring the format of row 15 (bordered/bold
#    label style in column A) and the shared-string label in column B.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.003347938170875
$ws.Range("D16").Value = 0.9733256173515474
$ws.Range("E16").Value = 1.001125728766453
$ws.Range("F16").Value = 1.003347938170875
$ws.Range("G16").Value = 0.9836962612942947
$ws.Range("H16").Value = 1.012687467361678
$ws.Range("I16").Value = 0.9998331504704114
$ws.Range("J16").Value = 0.9733256173515474
$ws.Range("K16").Value = 0.9872256730590003
$ws.Range("L16").Value = 0.9952868056149375
$ws.Range("M16").Value = 0.9956693605692098
